# Update the cryptocurrency price/volume table with the latest scraped
# values (GitHub Actions scheduled refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) and Volume(1h) (column E) updates -------------------
# Row => (Price, Volume%)   (Volume% left out where it did not change)

$ws.Range("D2").Value  = "24.771.56"
$ws.Range("E2").Value  = "  +1.64%  "

$ws.Range("D3").Value  = "1.701.59"
$ws.Range("E3").Value  = "  +1.54%  "

$ws.Range("D4").Value  = "1.008"
$ws.Range("E4").Value  = "  +0.35%  "

$ws.Range("D5").Value  = "311.12"
$ws.Range("E5").Value  = "  +1.83%  "

$ws.Range("D6").Value  = "1.003"
$ws.Range("E6").Value  = "  +0.29%  "

$ws.Range("D7").Value  = "0.3722"
$ws.Range("E7").Value  = "  +0.98%  "

$ws.Range("D8").Value  = "49.18"

$ws.Range("D9").Value  = "0.3418"
$ws.Range("E9").Value  = "  +0.06%  "

$ws.Range("D10").Value = "1.210"
$ws.Range("E10").Value = "  +4.67%  "

$ws.Range("D11").Value = "0.07482"
$ws.Range("E11").Value = "  +3.83%  "

$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("D13").Value = "20.94"
$ws.Range("E13").Value = "  +4.42%  "

$ws.Range("D14").Value = "6.298"
$ws.Range("E14").Value = "  +2.77%  "

$ws.Range("D15").Value = "7.009"
$ws.Range("E15").Value = "  +4.35%  "

$ws.Range("D16").Value = "1.700.35"
$ws.Range("E16").Value = "  +1.50%  "

$ws.Range("D17").Value = "0.00001125"
$ws.Range("E17").Value = "  +2.29%  "

$ws.Range("D18").Value = "0.06717"
$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").Value = "83.55"
$ws.Range("E20").Value = "  +4.11%  "

$ws.Range("D21").Value = "17.15"
$ws.Range("E21").Value = "  +3.89%  "

$ws.Range("D22").Value = "6.332"
$ws.Range("E22").Value = "  +4.02%  "

$ws.Range("D23").Value = "12.98"
$ws.Range("E23").Value = "  +6.94%  "

$ws.Range("D24").Value = "24.791.75"
$ws.Range("E24").Value = "  +1.93%  "

$ws.Range("D25").Value = "2.451"
$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("D26").Value = "2.763"
$ws.Range("E26").Value = "  +4.33%  "

$ws.Range("D27").Value = "20.24"
$ws.Range("E27").Value = "  +4.36%  "

$ws.Range("D28").Value = "148.95"
$ws.Range("E28").Value = "  -2.10%  "

$ws.Range("D29").Value = "131.74"
$ws.Range("E29").Value = "  +3.40%  "

# --- Rows 30 & 31 swap places (ranking order changed) ----------------------
# Row 30 was ImmutableX, Row 31 was WrappedliquidstakedEther2.0.
# They now swap order: WrappedliquidstakedEther2.0 moves to row 30,
# ImmutableX moves to row 31. The rank numbers in column A stay put.

$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.897.52"
$ws.Range("E30").Value = "  +2.01%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "1.244"
$ws.Range("E31").Value = "  +28.63%  "

$ws.Range("D32").Value = "6.724"
$ws.Range("E32").Value = "  +7.34%  "

$ws.Range("D33").Value = "4.230"
$ws.Range("E33").Value = "  +4.60%  "

$ws.Range("D34").Value = "13.64"
$ws.Range("E34").Value = "  +10.34%  "

$ws.Range("D35").Value = "1.765"
$ws.Range("E35").Value = "  +5.30%  "

$ws.Range("D36").Value = "0.08724"
$ws.Range("E36").Value = "  +3.35%  "

$ws.Range("D37").Value = "5.561"
$ws.Range("E37").Value = "  +4.81%  "

$ws.Range("D38").Value = "0.06617"
$ws.Range("E38").Value = "  +3.59%  "

$ws.Range("D39").Value = "9.054"
$ws.Range("E39").Value = "  +4.58%  "

$ws.Range("D40").Value = "0.02396"
$ws.Range("E40").Value = "  +3.68%  "

$ws.Range("D41").Value = "0.2220"
$ws.Range("E41").Value = "  +6.62%  "

$ws.Range("D42").Value = "1.271"
$ws.Range("E42").Value = "  +3.51%  "

$ws.Range("D43").Value = "0.6399"
$ws.Range("E43").Value = "  +5.53%  "

$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").Value = "13.79"
$ws.Range("E45").Value = "  +6.90%  "

$ws.Range("D46").Value = "0.6102"
$ws.Range("E46").Value = "  +4.25%  "

$ws.Range("D47").Value = "3.825"
$ws.Range("E47").Value = "  +1.95%  "

$ws.Range("D48").Value = "2.099"
$ws.Range("E48").Value = "  +4.22%  "

$ws.Range("D49").Value = "128.93"
$ws.Range("E49").Value = "  +2.48%  "

$ws.Range("D50").Value = "0.07276"
$ws.Range("E50").Value = "  +1.91%  "

$ws.Range("D51").Value = "79.39"
$ws.Range("E51").Value = "  +4.78%  "
